$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 34, pushing the existing data (old rows 34-57)
# down to rows 36-59.
$ws.Rows.Item(34).Insert()
$ws.Rows.Item(34).Insert()

# New row 34: weekly "Primera" quality record dated 2023-10-17 (serial 45216)
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 45216
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100112044
$ws.Range("G34").Value = "Perejil"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 150
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = 1500
$ws.Range("N34").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 750
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = "Hortaliza"

# New row 35: weekly "Segunda" quality record, same date
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C35").Value = "Arica y Parinacota"
$ws.Range("D35").Value = 45216
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 100112044
$ws.Range("G35").Value = "Perejil"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 220
$ws.Range("K35").Value = 1400
$ws.Range("L35").Value = 1400
$ws.Range("M35").Value = 1400
$ws.Range("N35").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 700
$ws.Range("Q35").Value = 2
$ws.Range("R35").Value = "Hortaliza"
